$d = $word.ActiveDocument

# 1. Rename the class/constructor references: "BuildingFinder" -> "JSCompiler"
#    (two occurrences, in the NAME and SYNOPSIS sections of the first man page)
$d.Content.Find.Execute("BuildingFinder", $true, $false, $false, $false, $false,
                         $true, 1, $false, "JSCompiler", 2) | Out-Null

# 2. The DESCRIPTION paragraph of the "run" section used to have its text split
#    into two runs around a "_GoBack" bookmark (". At " | bookmark | "first, parse...").
#    Re-join that text into a single run (this also removes the old bookmark that
#    sat between the two runs).
$oldSplitText = ". At first, parse the input source code with ANTLR4 and generate AST. And the type inference module will infer the type of each variables and functions. At last, generate LLVM IR with "
$d.Content.Find.Execute($oldSplitText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $oldSplitText, 2) | Out-Null

# 3. Re-create the "_GoBack" bookmark as a zero-length bookmark right at the end
#    of the DESCRIPTION paragraph of the first "JSCompiler" man page (after the
#    closing ".", before the paragraph mark).
#
#    Directly collapsing a Range to a position that sits exactly on a paragraph's
#    trailing paragraph-mark confuses Bookmarks.Add in this host, so we briefly
#    insert a one-character placeholder after the final ".", anchor the bookmark
#    just before that placeholder (now a perfectly ordinary mid-paragraph spot),
#    and then delete the placeholder again.
$descPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*output.ll*") {
        $descPara = $p
        break
    }
}

$endPos = $descPara.Range.End - 1
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")

$bkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bkRange) | Out-Null

$d.Range($endPos, $endPos + 1).Delete() | Out-Null
